$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos price/volume data
$ws.Range("D2").Value = "20.110.94"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.423.59"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9968"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.74"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3711"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3151"
$ws.Range("E8").Value = "  +2.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.74"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.063"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06560"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9973"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.548"
$ws.Range("E13").Value = "  +3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.17"
$ws.Range("E14").Value = "  +5.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.212"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "1.423.20"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05712"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9973"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.91"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.616"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.87"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.11"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.227"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").Value = "20.131.72"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.295"
$ws.Range("E26").Value = "  +4.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.45"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.35"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "1.583.25"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.969"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.303"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8288"
$ws.Range("E33").Value = "  -7.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07798"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("E35").Value = "  +11.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.921"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05866"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.912"
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9965"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.75"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02065"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1876"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5355"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.556"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.35"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.52"
$ws.Range("E47").Value = "  +5.77%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.792"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.042"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9988"
$ws.Range("E51").Value = "  -0.44%  "
